$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("B14")
$r.Font.Bold = $true
$r.Font.Size = 12
$r.VerticalAlignment = -4108

$r.Copy()
$dest = $ws.Range("B15:B17")
$dest.PasteSpecial(-4122)
Write-Host "done"
